$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 196, shifting existing rows 196-235 down to 197-236
$ws.Rows.Item(196).Insert()

# Row 196: new record inserted
$ws.Cells.Item(196, 1).Value2 = 9
$ws.Cells.Item(196, 2).Value2 = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(196, 3).Value2 = "Metropolitana"
$ws.Cells.Item(196, 4).Value2 = 44504
$ws.Cells.Item(196, 5).Value2 = 13
$ws.Cells.Item(196, 6).Value2 = 100112044
$ws.Cells.Item(196, 7).Value2 = "Perejil"
$ws.Cells.Item(196, 8).Value2 = "Sin especificar"
$ws.Cells.Item(196, 9).Value2 = "Primera"
$ws.Cells.Item(196, 10).Value2 = 144
$ws.Cells.Item(196, 11).Value2 = 8000
$ws.Cells.Item(196, 12).Value2 = 10000
$ws.Cells.Item(196, 13).Value2 = 9000
$ws.Cells.Item(196, 14).Value2 = "`$/docena de atados"
$ws.Cells.Item(196, 15).Value2 = "Región Metropolitana"
$ws.Cells.Item(196, 16).Value2 = 3000
$ws.Cells.Item(196, 17).Value2 = 3
$ws.Cells.Item(196, 18).Value2 = "Hortaliza"

# Rows 197-236: shifted values (only changed columns)
# Row 197
$ws.Cells.Item(197, 4).Value2 = 44384
$ws.Cells.Item(197, 10).Value2 = 34
$ws.Cells.Item(197, 11).Value2 = 7000
$ws.Cells.Item(197, 12).Value2 = 7000
$ws.Cells.Item(197, 13).Value2 = 7000
$ws.Cells.Item(197, 14).Value2 = "`$/caja 36 atados"
$ws.Cells.Item(197, 15).Value2 = "Región Metropolitana"
$ws.Cells.Item(197, 16).Value2 = 194
$ws.Cells.Item(197, 17).Value2 = 36

# Row 198
$ws.Cells.Item(198, 4).Value2 = 44384
$ws.Cells.Item(198, 10).Value2 = 70
$ws.Cells.Item(198, 11).Value2 = 15000
$ws.Cells.Item(198, 12).Value2 = 16000
$ws.Cells.Item(198, 13).Value2 = 15500
$ws.Cells.Item(198, 14).Value2 = "`$/docena de atados"
$ws.Cells.Item(198, 15).Value2 = "Región Metropolitana"
$ws.Cells.Item(198, 16).Value2 = 5167
$ws.Cells.Item(198, 17).Value2 = 3

# Row 199
$ws.Cells.Item(199, 4).Value2 = 44246
$ws.Cells.Item(199, 10).Value2 = 70
$ws.Cells.Item(199, 11).Value2 = 10000
$ws.Cells.Item(199, 12).Value2 = 10000
$ws.Cells.Item(199, 13).Value2 = 10000
$ws.Cells.Item(199, 14).Value2 = "`$/docena de atados"
$ws.Cells.Item(199, 15).Value2 = "Región Metropolitana"
$ws.Cells.Item(199, 16).Value2 = 3333
$ws.Cells.Item(199, 17).Value2 = 3

# Row 200
$ws.Cells.Item(200, 4).Value2 = 44487
$ws.Cells.Item(200, 10).Value2 = 61
$ws.Cells.Item(200, 11).Value2 = 10000
$ws.Cells.Item(200, 12).Value2 = 12000
$ws.Cells.Item(200, 13).Value2 = 10984
$ws.Cells.Item(200, 14).Value2 = "`$/docena de atados"
$ws.Cells.Item(200, 15).Value2 = "Región Metropolitana"
$ws.Cells.Item(200, 16).Value2 = 3661
$ws.Cells.Item(200, 17).Value2 = 3

# Row 201
$ws.Cells.Item(201, 4).Value2 = 44425
$ws.Cells.Item(201, 10).Value2 = 34
$ws.Cells.Item(201, 11).Value2 = 8000
$ws.Cells.Item(201, 12).Value2 = 8000
$ws.Cells.Item(201, 13).Value2 = 8000
$ws.Cells.Item(201, 14).Value2 = "`$/caja 36 atados"
$ws.Cells.Item(201, 15).Value2 = "Región Metropolitana"
$ws.Cells.Item(201, 16).Value2 = 222
$ws.Cells.Item(201, 17).Value2 = 36

# Row 202
$ws.Cells.Item(202, 4).Value2 = 44425
$ws.Cells.Item(202, 10).Value2 = 133
$ws.Cells.Item(202, 11).Value2 = 13000
$ws.Cells.Item(202, 12).Value2 = 15000
$ws.Cells.Item(202, 13).Value2 = 13992
$ws.Cells.Item(202, 14).Value2 = "`$/docena de atados"
$ws.Cells.Item(202, 15).Value2 = "Región Metropolitana"
$ws.Cells.Item(202, 16).Value2 = 4664
$ws.Cells.Item(202, 17).Value2 = 3

# Row 203
$ws.Cells.Item(203, 4).Value2 = 44343
$ws.Cells.Item(203, 10).Value2 = 40
$ws.Cells.Item(203, 11).Value2 = 6000
$ws.Cells.Item(203, 12).Value2 = 6000
$ws.Cells.Item(203, 13).Value2 = 6000
$ws.Cells.Item(203, 14).Value2 = "`$/caja 36 atados"
$ws.Cells.Item(203, 15).Value2 = "Región Metropolitana"
$ws.Cells.Item(203, 16).Value2 = 167
$ws.Cells.Item(203, 17).Value2 = 36

# Row 204
$ws.Cells.Item(204, 4).Value2 = 44343
$ws.Cells.Item(204, 10).Value2 = 70
$ws.Cells.Item(204, 11).Value2 = 11000
$ws.Cells.Item(204, 12).Value2 = 12000
$ws.Cells.Item(204, 13).Value2 = 11500
$ws.Cells.Item(204, 14).Value2 = "`$/docena de atados"
$ws.Cells.Item(204, 15).Value2 = "Región Metropolitana"
$ws.Cells.Item(204, 16).Value2 = 3833
$ws.Cells.Item(204, 17).Value2 = 3

# Row 205
$ws.Cells.Item(205, 4).Value2 = 44370
$ws.Cells.Item(205, 10).Value2 = 34
$ws.Cells.Item(205, 11).Value2 = 6000
$ws.Cells.Item(205, 12).Value2 = 6000
$ws.Cells.Item(205, 13).Value2 = 6000
$ws.Cells.Item(205, 14).Value2 = "`$/caja 36 atados"
$ws.Cells.Item(205, 15).Value2 = "Región Metropolitana"
$ws.Cells.Item(205, 16).Value2 = 167
$ws.Cells.Item(205, 17).Value2 = 36

# Row 206
$ws.Cells.Item(206, 4).Value2 = 44370
$ws.Cells.Item(206, 10).Value2 = 52
$ws.Cells.Item(206, 11).Value2 = 11000
$ws.Cells.Item(206, 12).Value2 = 11500
$ws.Cells.Item(206, 13).Value2 = 11250
$ws.Cells.Item(206, 14).Value2 = "`$/docena de atados"
$ws.Cells.Item(206, 15).Value2 = "Región Metropolitana"
$ws.Cells.Item(206, 16).Value2 = 3750
$ws.Cells.Item(206, 17).Value2 = 3

# Row 207
$ws.Cells.Item(207, 4).Value2 = 44449
$ws.Cells.Item(207, 10).Value2 = 106
$ws.Cells.Item(207, 11).Value2 = 9000
$ws.Cells.Item(207, 12).Value2 = 10000
$ws.Cells.Item(207, 13).Value2 = 9500
$ws.Cells.Item(207, 14).Value2 = "`$/docena de atados"
$ws.Cells.Item(207, 15).Value2 = "Región Metropolitana"
$ws.Cells.Item(207, 16).Value2 = 3167
$ws.Cells.Item(207, 17).Value2 = 3

# Row 208
$ws.Cells.Item(208, 4).Value2 = 44168
$ws.Cells.Item(208, 10).Value2 = 80
$ws.Cells.Item(208, 11).Value2 = 8000
$ws.Cells.Item(208, 12).Value2 = 10000
$ws.Cells.Item(208, 13).Value2 = 9250
$ws.Cells.Item(208, 14).Value2 = "`$/docena de atados"
$ws.Cells.Item(208, 15).Value2 = "Región Metropolitana"
$ws.Cells.Item(208, 16).Value2 = 3083
$ws.Cells.Item(208, 17).Value2 = 3

# Row 209
$ws.Cells.Item(209, 4).Value2 = 44175
$ws.Cells.Item(209, 10).Value2 = 80
$ws.Cells.Item(209, 11).Value2 = 9000
$ws.Cells.Item(209, 12).Value2 = 10000
$ws.Cells.Item(209, 13).Value2 = 9500
$ws.Cells.Item(209, 14).Value2 = "`$/docena de atados"
$ws.Cells.Item(209, 15).Value2 = "Región Metropolitana"
$ws.Cells.Item(209, 16).Value2 = 3167
$ws.Cells.Item(209, 17).Value2 = 3

# Row 210
$ws.Cells.Item(210, 4).Value2 = 44392
$ws.Cells.Item(210, 10).Value2 = 52
$ws.Cells.Item(210, 11).Value2 = 7000
$ws.Cells.Item(210, 12).Value2 = 7000
$ws.Cells.Item(210, 13).Value2 = 7000
$ws.Cells.Item(210, 14).Value2 = "`$/caja 36 atados"
$ws.Cells.Item(210, 15).Value2 = "Región Metropolitana"
$ws.Cells.Item(210, 16).Value2 = 194
$ws.Cells.Item(210, 17).Value2 = 36

# Row 211
$ws.Cells.Item(211, 4).Value2 = 44392
$ws.Cells.Item(211, 10).Value2 = 70
$ws.Cells.Item(211, 11).Value2 = 12000
$ws.Cells.Item(211, 12).Value2 = 13000
$ws.Cells.Item(211, 13).Value2 = 12500
$ws.Cells.Item(211, 14).Value2 = "`$/docena de atados"
$ws.Cells.Item(211, 15).Value2 = "Región Metropolitana"
$ws.Cells.Item(211, 16).Value2 = 4167
$ws.Cells.Item(211, 17).Value2 = 3

# Row 212
$ws.Cells.Item(212, 4).Value2 = 44286
$ws.Cells.Item(212, 10).Value2 = 70
$ws.Cells.Item(212, 11).Value2 = 12000
$ws.Cells.Item(212, 12).Value2 = 12000
$ws.Cells.Item(212, 13).Value2 = 12000
$ws.Cells.Item(212, 14).Value2 = "`$/docena de atados"
$ws.Cells.Item(212, 15).Value2 = "Región Metropolitana"
$ws.Cells.Item(212, 16).Value2 = 4000
$ws.Cells.Item(212, 17).Value2 = 3

# Row 213
$ws.Cells.Item(213, 4).Value2 = 44400
$ws.Cells.Item(213, 10).Value2 = 34
$ws.Cells.Item(213, 11).Value2 = 6500
$ws.Cells.Item(213, 12).Value2 = 7000
$ws.Cells.Item(213, 13).Value2 = 6750
$ws.Cells.Item(213, 14).Value2 = "`$/caja 36 atados"
$ws.Cells.Item(213, 15).Value2 = "Región Metropolitana"
$ws.Cells.Item(213, 16).Value2 = 188
$ws.Cells.Item(213, 17).Value2 = 36

# Row 214
$ws.Cells.Item(214, 4).Value2 = 44400
$ws.Cells.Item(214, 10).Value2 = 70
$ws.Cells.Item(214, 11).Value2 = 13000
$ws.Cells.Item(214, 12).Value2 = 14000
$ws.Cells.Item(214, 13).Value2 = 13500
$ws.Cells.Item(214, 14).Value2 = "`$/docena de atados"
$ws.Cells.Item(214, 15).Value2 = "Región Metropolitana"
$ws.Cells.Item(214, 16).Value2 = 4500
$ws.Cells.Item(214, 17).Value2 = 3

# Row 215
$ws.Cells.Item(215, 4).Value2 = 44484
$ws.Cells.Item(215, 10).Value2 = 97
$ws.Cells.Item(215, 11).Value2 = 8000
$ws.Cells.Item(215, 12).Value2 = 10000
$ws.Cells.Item(215, 13).Value2 = 8990
$ws.Cells.Item(215, 14).Value2 = "`$/docena de atados"
$ws.Cells.Item(215, 15).Value2 = "Región Metropolitana"
$ws.Cells.Item(215, 16).Value2 = 2997
$ws.Cells.Item(215, 17).Value2 = 3

# Row 216
$ws.Cells.Item(216, 4).Value2 = 44181
$ws.Cells.Item(216, 10).Value2 = 70
$ws.Cells.Item(216, 11).Value2 = 8000
$ws.Cells.Item(216, 12).Value2 = 10000
$ws.Cells.Item(216, 13).Value2 = 8857
$ws.Cells.Item(216, 14).Value2 = "`$/docena de atados"
$ws.Cells.Item(216, 15).Value2 = "Región Metropolitana"
$ws.Cells.Item(216, 16).Value2 = 2952
$ws.Cells.Item(216, 17).Value2 = 3

# Row 217
$ws.Cells.Item(217, 4).Value2 = 44494
$ws.Cells.Item(217, 10).Value2 = 80
$ws.Cells.Item(217, 11).Value2 = 10000
$ws.Cells.Item(217, 12).Value2 = 10000
$ws.Cells.Item(217, 13).Value2 = 10000
$ws.Cells.Item(217, 14).Value2 = "`$/docena de atados"
$ws.Cells.Item(217, 15).Value2 = "Provincia de Chacabuco"
$ws.Cells.Item(217, 16).Value2 = 3333
$ws.Cells.Item(217, 17).Value2 = 3

# Row 218
$ws.Cells.Item(218, 4).Value2 = 44342
$ws.Cells.Item(218, 10).Value2 = 30
$ws.Cells.Item(218, 11).Value2 = 6000
$ws.Cells.Item(218, 12).Value2 = 6000
$ws.Cells.Item(218, 13).Value2 = 6000
$ws.Cells.Item(218, 14).Value2 = "`$/caja 36 atados"
$ws.Cells.Item(218, 15).Value2 = "Región Metropolitana"
$ws.Cells.Item(218, 16).Value2 = 167
$ws.Cells.Item(218, 17).Value2 = 36

# Row 219
$ws.Cells.Item(219, 4).Value2 = 44342
$ws.Cells.Item(219, 10).Value2 = 80
$ws.Cells.Item(219, 11).Value2 = 11000
$ws.Cells.Item(219, 12).Value2 = 12000
$ws.Cells.Item(219, 13).Value2 = 11500
$ws.Cells.Item(219, 14).Value2 = "`$/docena de atados"
$ws.Cells.Item(219, 15).Value2 = "Región Metropolitana"
$ws.Cells.Item(219, 16).Value2 = 3833
$ws.Cells.Item(219, 17).Value2 = 3

# Row 220
$ws.Cells.Item(220, 4).Value2 = 44328
$ws.Cells.Item(220, 10).Value2 = 52
$ws.Cells.Item(220, 11).Value2 = 12000
$ws.Cells.Item(220, 12).Value2 = 12000
$ws.Cells.Item(220, 13).Value2 = 12000
$ws.Cells.Item(220, 14).Value2 = "`$/docena de atados"
$ws.Cells.Item(220, 15).Value2 = "Región Metropolitana"
$ws.Cells.Item(220, 16).Value2 = 4000
$ws.Cells.Item(220, 17).Value2 = 3

# Row 221
$ws.Cells.Item(221, 4).Value2 = 44301
$ws.Cells.Item(221, 10).Value2 = 70
$ws.Cells.Item(221, 11).Value2 = 12000
$ws.Cells.Item(221, 12).Value2 = 12000
$ws.Cells.Item(221, 13).Value2 = 12000
$ws.Cells.Item(221, 14).Value2 = "`$/docena de atados"
$ws.Cells.Item(221, 15).Value2 = "Región Metropolitana"
$ws.Cells.Item(221, 16).Value2 = 4000
$ws.Cells.Item(221, 17).Value2 = 3

# Row 222
$ws.Cells.Item(222, 4).Value2 = 44330
$ws.Cells.Item(222, 10).Value2 = 25
$ws.Cells.Item(222, 11).Value2 = 6000
$ws.Cells.Item(222, 12).Value2 = 6000
$ws.Cells.Item(222, 13).Value2 = 6000
$ws.Cells.Item(222, 14).Value2 = "`$/caja 36 atados"
$ws.Cells.Item(222, 15).Value2 = "Región Metropolitana"
$ws.Cells.Item(222, 16).Value2 = 167
$ws.Cells.Item(222, 17).Value2 = 36

# Row 223
$ws.Cells.Item(223, 4).Value2 = 44330
$ws.Cells.Item(223, 10).Value2 = 50
$ws.Cells.Item(223, 11).Value2 = 12000
$ws.Cells.Item(223, 12).Value2 = 12000
$ws.Cells.Item(223, 13).Value2 = 12000
$ws.Cells.Item(223, 14).Value2 = "`$/docena de atados"
$ws.Cells.Item(223, 15).Value2 = "Región Metropolitana"
$ws.Cells.Item(223, 16).Value2 = 4000
$ws.Cells.Item(223, 17).Value2 = 3

# Row 224
$ws.Cells.Item(224, 4).Value2 = 44270
$ws.Cells.Item(224, 10).Value2 = 60
$ws.Cells.Item(224, 11).Value2 = 14000
$ws.Cells.Item(224, 12).Value2 = 14000
$ws.Cells.Item(224, 13).Value2 = 14000
$ws.Cells.Item(224, 14).Value2 = "`$/docena de atados"
$ws.Cells.Item(224, 15).Value2 = "Región Metropolitana"
$ws.Cells.Item(224, 16).Value2 = 4667
$ws.Cells.Item(224, 17).Value2 = 3

# Row 225
$ws.Cells.Item(225, 4).Value2 = 44295
$ws.Cells.Item(225, 10).Value2 = 50
$ws.Cells.Item(225, 11).Value2 = 13000
$ws.Cells.Item(225, 12).Value2 = 13000
$ws.Cells.Item(225, 13).Value2 = 13000
$ws.Cells.Item(225, 14).Value2 = "`$/docena de atados"
$ws.Cells.Item(225, 15).Value2 = "Región Metropolitana"
$ws.Cells.Item(225, 16).Value2 = 4333
$ws.Cells.Item(225, 17).Value2 = 3

# Row 226
$ws.Cells.Item(226, 4).Value2 = 44217
$ws.Cells.Item(226, 10).Value2 = 95
$ws.Cells.Item(226, 11).Value2 = 8000
$ws.Cells.Item(226, 12).Value2 = 9000
$ws.Cells.Item(226, 13).Value2 = 8368
$ws.Cells.Item(226, 14).Value2 = "`$/docena de atados"
$ws.Cells.Item(226, 15).Value2 = "Región Metropolitana"
$ws.Cells.Item(226, 16).Value2 = 2789
$ws.Cells.Item(226, 17).Value2 = 3

# Row 227
$ws.Cells.Item(227, 4).Value2 = 44421
$ws.Cells.Item(227, 10).Value2 = 133
$ws.Cells.Item(227, 11).Value2 = 10000
$ws.Cells.Item(227, 12).Value2 = 11000
$ws.Cells.Item(227, 13).Value2 = 10496
$ws.Cells.Item(227, 14).Value2 = "`$/docena de atados"
$ws.Cells.Item(227, 15).Value2 = "Región Metropolitana"
$ws.Cells.Item(227, 16).Value2 = 3499
$ws.Cells.Item(227, 17).Value2 = 3

# Row 228
$ws.Cells.Item(228, 4).Value2 = 44383
$ws.Cells.Item(228, 10).Value2 = 34
$ws.Cells.Item(228, 11).Value2 = 6000
$ws.Cells.Item(228, 12).Value2 = 6000
$ws.Cells.Item(228, 13).Value2 = 6000
$ws.Cells.Item(228, 14).Value2 = "`$/caja 36 atados"
$ws.Cells.Item(228, 15).Value2 = "Región Metropolitana"
$ws.Cells.Item(228, 16).Value2 = 167
$ws.Cells.Item(228, 17).Value2 = 36

# Row 229
$ws.Cells.Item(229, 4).Value2 = 44383
$ws.Cells.Item(229, 10).Value2 = 52
$ws.Cells.Item(229, 11).Value2 = 11000
$ws.Cells.Item(229, 12).Value2 = 11500
$ws.Cells.Item(229, 13).Value2 = 11250
$ws.Cells.Item(229, 14).Value2 = "`$/docena de atados"
$ws.Cells.Item(229, 15).Value2 = "Región Metropolitana"
$ws.Cells.Item(229, 16).Value2 = 3750
$ws.Cells.Item(229, 17).Value2 = 3

# Row 230
$ws.Cells.Item(230, 4).Value2 = 44244
$ws.Cells.Item(230, 10).Value2 = 50
$ws.Cells.Item(230, 11).Value2 = 10000
$ws.Cells.Item(230, 12).Value2 = 10000
$ws.Cells.Item(230, 13).Value2 = 10000
$ws.Cells.Item(230, 14).Value2 = "`$/docena de atados"
$ws.Cells.Item(230, 15).Value2 = "Región Metropolitana"
$ws.Cells.Item(230, 16).Value2 = 3333
$ws.Cells.Item(230, 17).Value2 = 3

# Row 231
$ws.Cells.Item(231, 4).Value2 = 44307
$ws.Cells.Item(231, 10).Value2 = 70
$ws.Cells.Item(231, 11).Value2 = 12000
$ws.Cells.Item(231, 12).Value2 = 12000
$ws.Cells.Item(231, 13).Value2 = 12000
$ws.Cells.Item(231, 14).Value2 = "`$/docena de atados"
$ws.Cells.Item(231, 15).Value2 = "Región Metropolitana"
$ws.Cells.Item(231, 16).Value2 = 4000
$ws.Cells.Item(231, 17).Value2 = 3

# Row 232
$ws.Cells.Item(232, 4).Value2 = 44273
$ws.Cells.Item(232, 10).Value2 = 70
$ws.Cells.Item(232, 11).Value2 = 10000
$ws.Cells.Item(232, 12).Value2 = 10000
$ws.Cells.Item(232, 13).Value2 = 10000
$ws.Cells.Item(232, 14).Value2 = "`$/docena de atados"
$ws.Cells.Item(232, 15).Value2 = "Región Metropolitana"
$ws.Cells.Item(232, 16).Value2 = 3333
$ws.Cells.Item(232, 17).Value2 = 3

# Row 233
$ws.Cells.Item(233, 4).Value2 = 44433
$ws.Cells.Item(233, 10).Value2 = 97
$ws.Cells.Item(233, 11).Value2 = 11000
$ws.Cells.Item(233, 12).Value2 = 12000
$ws.Cells.Item(233, 13).Value2 = 11495
$ws.Cells.Item(233, 14).Value2 = "`$/docena de atados"
$ws.Cells.Item(233, 15).Value2 = "Región Metropolitana"
$ws.Cells.Item(233, 16).Value2 = 3832
$ws.Cells.Item(233, 17).Value2 = 3

# Row 234
$ws.Cells.Item(234, 4).Value2 = 44302
$ws.Cells.Item(234, 10).Value2 = 50
$ws.Cells.Item(234, 11).Value2 = 12000
$ws.Cells.Item(234, 12).Value2 = 12000
$ws.Cells.Item(234, 13).Value2 = 12000
$ws.Cells.Item(234, 14).Value2 = "`$/docena de atados"
$ws.Cells.Item(234, 15).Value2 = "Región Metropolitana"
$ws.Cells.Item(234, 16).Value2 = 4000
$ws.Cells.Item(234, 17).Value2 = 3

# Row 235
$ws.Cells.Item(235, 4).Value2 = 44179
$ws.Cells.Item(235, 10).Value2 = 35
$ws.Cells.Item(235, 11).Value2 = 8000
$ws.Cells.Item(235, 12).Value2 = 9000
$ws.Cells.Item(235, 13).Value2 = 8429
$ws.Cells.Item(235, 14).Value2 = "`$/docena de atados"
$ws.Cells.Item(235, 15).Value2 = "Región Metropolitana"
$ws.Cells.Item(235, 16).Value2 = 2810
$ws.Cells.Item(235, 17).Value2 = 3

# Row 236
$ws.Cells.Item(236, 4).Value2 = 44491
$ws.Cells.Item(236, 10).Value2 = 106
$ws.Cells.Item(236, 11).Value2 = 9000
$ws.Cells.Item(236, 12).Value2 = 10000
$ws.Cells.Item(236, 13).Value2 = 9500
$ws.Cells.Item(236, 14).Value2 = "`$/docena de atados"
$ws.Cells.Item(236, 15).Value2 = "Región Metropolitana"
$ws.Cells.Item(236, 16).Value2 = 3167
$ws.Cells.Item(236, 17).Value2 = 3
